$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 19 (the "Arduino Experiment 4X4
# Matrix" row), shifting it - and everything below it - down by one. Excel's
# default Insert copies formatting down from the row above (row 18), which
# already matches the desired styling for the new row's A/B cells, and the
# totals formulas (SUM / division) get their ranges auto-adjusted for the
# shift.
$ws.Rows.Item(19).Insert()

# The "Date" column stores entries such as "26.09.2022" as plain text (not
# real dates). Assigning a day/month-ambiguous string like "11.10.2022"
# straight to .Value gets auto-parsed into a date serial by the automation
# layer, which would also fabricate a brand-new number-format style. Route
# it through a throwaway TRIM() formula and a values-only paste instead, so
# the literal text lands in the cell without disturbing its (already
# correct) style or creating any new style entries.
$ws.Range("Z1").Formula = '=TRIM("11.10.2022 ")'
$ws.Range("Z1").Copy()
$ws.Range("A19").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("B19").Value = 0.79166666666666663
$ws.Range("C19").Value = "Research Keyboards"
$ws.Range("D19").Value = "Research"
$ws.Range("E19").Value = 80
$ws.Range("F19").Value = "Research on Keyboard Scanning and Create Flowchart"

# Match the author's final selection state.
$ws.Range("F20").Select()
